$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "495÷2="; New = "200÷8=" },
    @{ Old = "507÷4="; New = "173÷8=" },
    @{ Old = "496÷2="; New = "775÷5=" },
    @{ Old = "698÷4="; New = "997÷5=" },
    @{ Old = "590÷2="; New = "509÷4=" },
    @{ Old = "707÷3="; New = "186÷3=" },
    @{ Old = "658÷7="; New = "417÷4=" },
    @{ Old = "831÷8="; New = "175÷7=" },
    @{ Old = "111÷6="; New = "239÷2=" },
    @{ Old = "878÷4="; New = "258÷7=" },
    @{ Old = "825÷8="; New = "728÷5=" },
    @{ Old = "152÷9="; New = "711÷3=" },
    @{ Old = "887÷2="; New = "491÷5=" },
    @{ Old = "783÷2="; New = "732÷3=" },
    @{ Old = "590÷7="; New = "116÷7=" },
    @{ Old = "261÷2="; New = "393÷7=" },
    @{ Old = "332÷7="; New = "261÷5=" },
    @{ Old = "219÷2="; New = "389÷9=" },
    @{ Old = "195÷6="; New = "288÷3=" },
    @{ Old = "202÷4="; New = "905÷6=" },
    @{ Old = "464÷2="; New = "154÷5=" },
    @{ Old = "833÷9="; New = "314÷6=" },
    @{ Old = "403÷3="; New = "845÷4=" },
    @{ Old = "728÷2="; New = "208÷8=" },
    @{ Old = "526÷4="; New = "638÷5=" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.New, 2)
}
